$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Fitness) for rows 2 through 21 from 7586 to 7569
$ws.Range("C2:C21").Value = 7569
